$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 20 first, shifting old rows 20,21,... down to 21,22,...
$ws.Rows.Item(20).Insert()

# Row 15: give it a custom height (16.5)
$ws.Rows.Item(15).RowHeight = 16.5

# Row 18: change label text
$ws.Range("A18").Value = "R1 <- INTER"

# Row 19: new data pattern
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("T19").Value = 1

# Row 20: new label row "R1 <- R1 + R2" (mirrors the style of row 18/16 label rows)
$ws.Range("A20:B20").Merge()
$ws.Range("A20").Value = "R1 <- R1 + R2"
$ws.Range("A20:B20").Style = $ws.Range("A18:B18").Style
$ws.Range("H20").Style = $ws.Range("H18").Style

# Row 21: data pattern row
$ws.Range("A21").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "00000"
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 0

# Row 22: data pattern row (matches the original row20 pattern)
$ws.Range("A22").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = "00001"
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("T22").Value = 0

# Row 23: data pattern row (matches the original row21 pattern)
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = "00000"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 0
$ws.Range("S23").Value = 0
$ws.Range("T23").Value = 0

# New rows 35 and 36 appended at the end (empty placeholder rows, same style as row 34)
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"

# Update selection
$ws.Range("T18").Select()
